$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before the "总计" (summary) sheet.
#    Seed it by copying the "2021-Q4" sheet (same column layout / styling),
#    then overwrite the data cells with the 2022-Q1 numbers.
# ---------------------------------------------------------------------------
$totalSheetName = "总计"
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item($totalSheetName))
$newSheet.Name = "2022-Q1"

# Seed layout/styling from "2021-Q4" (same columns). Copy the header row and
# the data block separately so we don't drag along a spurious blank A1 cell.
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$q4Sheet.Range("A2:H5").Copy($newSheet.Range("A2:H5"))

# Fund-code column (B) and the decimal-looking text columns (D,E,F,G) must
# stay TEXT (they were stored as text in the source data), so they are
# written through Formula with a leading apostrophe to force text, exactly
# like the existing quarter sheets. Column C (fund name) is plain text
# already, and columns A/H are genuine numbers.

# Row 2: 006102 / 浙商丰利增强债券
$newSheet.Cells.Item(2,2).Formula = "'006102"
$newSheet.Cells.Item(2,3).Value2 = "浙商丰利增强债券"
$newSheet.Cells.Item(2,4).Formula = "'48.75"
$newSheet.Cells.Item(2,5).Formula = "'47.92"
$newSheet.Cells.Item(2,6).Formula = "'2.59"
$newSheet.Cells.Item(2,7).Formula = "'1.2626"
$newSheet.Cells.Item(2,8).Value2 = 4

# Row 3: 688888 / 浙商聚潮产业成长混合
$newSheet.Cells.Item(3,2).Formula = "'688888"
$newSheet.Cells.Item(3,3).Value2 = "浙商聚潮产业成长混合"
$newSheet.Cells.Item(3,4).Formula = "'8.25"
$newSheet.Cells.Item(3,5).Formula = "'93.40"
$newSheet.Cells.Item(3,6).Formula = "'6.03"
$newSheet.Cells.Item(3,7).Formula = "'0.4975"
$newSheet.Cells.Item(3,8).Value2 = 5

# Row 4: 010381 / 浙商智选价值混合A
$newSheet.Cells.Item(4,2).Formula = "'010381"
$newSheet.Cells.Item(4,3).Value2 = "浙商智选价值混合A"
$newSheet.Cells.Item(4,4).Formula = "'2.92"
$newSheet.Cells.Item(4,5).Formula = "'93.43"
$newSheet.Cells.Item(4,6).Formula = "'5.68"
$newSheet.Cells.Item(4,7).Formula = "'0.1659"
$newSheet.Cells.Item(4,8).Value2 = 5

# Row 5: 010382 / 浙商智选价值混合C
$newSheet.Cells.Item(5,2).Formula = "'010382"
$newSheet.Cells.Item(5,3).Value2 = "浙商智选价值混合C"
$newSheet.Cells.Item(5,4).Formula = "'0.34"
$newSheet.Cells.Item(5,5).Formula = "'93.43"
$newSheet.Cells.Item(5,6).Formula = "'5.68"
$newSheet.Cells.Item(5,7).Formula = "'0.0193"
$newSheet.Cells.Item(5,8).Value2 = 5

# ---------------------------------------------------------------------------
# 2. Add a 2022-Q1 summary row at the top of the "总计" sheet's data, pushing
#    the existing rows down by one (and bumping their row-index column A).
# ---------------------------------------------------------------------------

# Re-resolve "总计" by name now that the new sheet has been inserted - a
# reference captured before Worksheets.Add() tracks its numeric position, and
# Add() shifted "总计" from index 5 to index 6.
$totalSheet = $wb.Worksheets.Item($totalSheetName)

# Shift existing data rows 2-5 down to rows 3-6 (copy with Value2 so we don't
# depend on row-insert semantics; columns A/C/D are numbers, B is text).
for ($r = 5; $r -ge 2; $r--) {
    $dest = $r + 1
    $totalSheet.Cells.Item($dest, 1).Value2 = $r - 1
    $totalSheet.Cells.Item($dest, 2).Value2 = $totalSheet.Cells.Item($r, 2).Value2
    $totalSheet.Cells.Item($dest, 3).Value2 = $totalSheet.Cells.Item($r, 3).Value2
    $totalSheet.Cells.Item($dest, 4).Value2 = $totalSheet.Cells.Item($r, 4).Value2
}

# Row 6 (the old row 5) lost the bold/bordered index style when it was
# created fresh by the loop above - copy it over from row 5 (which still
# carries the correct style) so the whole A column stays consistently
# formatted.
$totalSheet.Range("A5").Copy($totalSheet.Range("A6"))
$totalSheet.Cells.Item(6, 1).Value2 = 4

# Write the new 2022-Q1 summary values into row 2.
$totalSheet.Cells.Item(2, 1).Value2 = 0
$totalSheet.Cells.Item(2, 2).Value2 = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value2 = 4
$totalSheet.Cells.Item(2, 4).Value2 = 1.95
